# Recolor the "cose da fare" highlighted notes from lightGray to cyan.
#
# The affected paragraphs are the block that starts with "Nel metodo
# sendData..." and ends with "...il loro avvio." Two of those paragraphs
# are "Titolo2" headings ("Scontrol" and "Riferimenti e configurazione")
# whose paragraph mark itself must also carry the cyan highlight, and one
# paragraph ("Puoi commentare...") previously had no <w:pPr> at all and
# needs one added (carrying the cyan highlight on the paragraph mark).
# The final paragraph additionally has its last two runs
# (". Lasci comunque ... avvio" + ".") merged into a single highlighted
# run ending in a period.

$d = $word.ActiveDocument

function Set-ParaBodyHighlight($paraIndex, $color) {
    # Recolor only the paragraph's text runs, leaving the paragraph mark
    # (and therefore any <w:pPr>) untouched.
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $body = $d.Range($full.Start, $full.End - 1)
    $body.Font.HighlightColorIndex = $color
}

function Set-ParaFullHighlight($paraIndex, $color) {
    # Recolor the paragraph's text runs AND the paragraph mark, which is
    # what materializes/updates the <w:rPr> inside <w:pPr>.
    $p = $d.Paragraphs($paraIndex)
    $p.Range.Font.HighlightColorIndex = $color
}

# "Nel metodo sendData, anziché mandare i dati tramite contact ..."
Set-ParaBodyHighlight 40 "cyan"

# "nel metodo MandaComando a Scontrol anziché il metodo userCmdDemand ..."
Set-ParaBodyHighlight 42 "cyan"

# "Scontrol" (Titolo2 heading) -- paragraph mark formatting also changes.
Set-ParaFullHighlight 43 "cyan"

# "Puoi commentare ciò che avveniva nei metodi Dojob e run del thread ..."
# -- this paragraph gains a brand-new <w:pPr> (it had none before).
Set-ParaFullHighlight 44 "cyan"

# "Riferimenti e configurazione" (Titolo2 heading) -- paragraph mark too.
Set-ParaFullHighlight 45 "cyan"

# "Sia i sensori che UserCmd devono avere ... e il loro avvio."
Set-ParaBodyHighlight 46 "cyan"

# Merge the trailing ". Lasci comunque invariato in Edi l'assegnazione di
# nomi ai processi e il loro avvio" run with the separate, unformatted
# "." run right after it, producing one run ending in "avvio." that
# carries the cyan highlight.
$p46 = $d.Paragraphs(46)
$searchRange = $p46.Range
$find = $searchRange.Find
$find.ClearFormatting()
$find.Text = ". Lasci comunque invariato in Edi l" + [char]0x2019 + "assegnazione di nomi ai processi e il loro avvio"
$found = $find.Execute()
if ($found) {
    $searchRange.InsertAfter(".")
    $leftover = $d.Range($searchRange.End, $searchRange.End + 1)
    if ($leftover.Text -eq ".") {
        $leftover.Delete()
    }
    $searchRange.Font.HighlightColorIndex = "cyan"
}
